# form_import_siswa.xlsx - rework the import template:
#  - add columns for NISN, NIK, Nama Ibu Kandung, Desa, Kabupaten, Kecamatan
#  - rename "No. KK/PKH/PIP" -> "Nomor KK/PKH/PIP"
#  - reorder all headers into the new column layout (A..O)
#  - give the new "Nomor PIP" column (E) the same text format + border as
#    the other "text id" columns (A-D)
#  - give every new column (K..O) the same border-only look as the other
#    data columns
#  - resize columns to fit their (new) header text
#  - move the active selection to L3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1) - new column order
# ---------------------------------------------------------------------
$headers = @{
    "A1" = "NISN"
    "B1" = "NIK"
    "C1" = "Nomor KK"
    "D1" = "Nomor PKH"
    "E1" = "Nomor PIP"
    "F1" = "Nama Siswa"
    "G1" = "Nama Ibu Kandung"
    "H1" = "Tempat Lahir"
    "I1" = "Tanggal Lahir"
    "J1" = "Alamat"
    "K1" = "Desa"
    "L1" = "Kabupaten"
    "M1" = "Kecamatan"
    "N1" = "Jenis Kelamin (L/P)"
    "O1" = "Status (WNI/WNA)"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value2 = $headers[$addr]
}

# Give the new K1:O1 headers the exact same look (Good cell style, thin
# border, centered) as the rest of the header row, by copying the
# formatting from an existing header cell - the constant xlPasteFormats
# (-4122) paste type copies formats only, values are left untouched.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1:O1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2. Data rows (2-11)
# ---------------------------------------------------------------------
# Column E ("Nomor PIP") now gets the same look as the other "id" style
# columns (A-D): thin border + text number format - copy it from A2:A11.
$ws.Range("A2:A11").Copy() | Out-Null
$ws.Range("E2:E11").PasteSpecial(-4122) | Out-Null

# New columns K:O (Desa, Kabupaten, Kecamatan, Jenis Kelamin, Status) get
# the plain bordered look used by F:J - copy it across.
$ws.Range("F2:J11").Copy() | Out-Null
$ws.Range("K2:O11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Column widths - fit the (new) header text, same as the sheet author
#    would get from "AutoFit Column Width". Columns A/B keep the sheet's
#    standard (default) width, untouched.
# ---------------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 8.92
$ws.Columns("D:D").ColumnWidth = 10.26
$ws.Columns("E:E").ColumnWidth = 9.59
$ws.Columns("F:F").ColumnWidth = 11.09
$ws.Columns("G:G").ColumnWidth = 17.42
$ws.Columns("H:H").ColumnWidth = 11.59
$ws.Columns("I:I").ColumnWidth = 11.92
$ws.Columns("L:L").ColumnWidth = 9.59
$ws.Columns("M:M").ColumnWidth = 10.09
$ws.Columns("N:N").ColumnWidth = 17.09
$ws.Columns("O:O").ColumnWidth = 16.42

# ---------------------------------------------------------------------
# 4. Selection - the author left the cursor on L3
# ---------------------------------------------------------------------
$ws.Range("L3").Select() | Out-Null

Write-Output "form_import_siswa template updated"
